$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) values: these are numeric-looking strings that must be
# preserved verbatim as text (they use "." as a thousands separator in some
# rows, and Excel would otherwise coerce "1.00" -> 1, "0.0458" -> 4.58E-02, etc).
# Temporarily force Text number format, assign the literal string, then restore
# the "Normal" style so no stray formatting is left behind on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.819.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.814.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '353.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.598'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.237.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.833.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.896'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.779.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0993'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '50.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0454'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '33.68'
$ws.Range("D33").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.90'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.26'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '127.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.079.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.941'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.51'
$ws.Range("D51").Style = "Normal"

# --- Coin name / link / Volume(1h) columns: plain text, safe to assign directly.
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +6.58%  '
$ws.Range("E6").Value = '  -1.83%  '
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +4.48%  '
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("E13").Value = '  +1.43%  '
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("E16").Value = '  +2.00%  '
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("E18").Value = '  +0.81%  '
$ws.Range("E19").Value = '  +8.67%  '
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("E22").Value = '  +2.18%  '
$ws.Range("E23").Value = '  -2.76%  '
$ws.Range("E24").Value = '  +0.50%  '
$ws.Range("E25").Value = '  +4.58%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("E29").Value = '  +1.82%  '
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("B31").Value = 'OKB'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("E31").Value = '  +1.63%  '
$ws.Range("B32").Value = 'VeChain'
$ws.Range("C32").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E32").Value = '  +30.78%  '
$ws.Range("E33").Value = '  -3.82%  '
$ws.Range("E34").Value = '  +5.01%  '
$ws.Range("E35").Value = '  +1.10%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("E40").Value = '  -4.00%  '
$ws.Range("E41").Value = '  +1.68%  '
$ws.Range("E42").Value = '  +5.18%  '
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("E44").Value = '  +1.48%  '
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E47").Value = '  +0.71%  '
$ws.Range("E49").Value = '  +2.97%  '
$ws.Range("E50").Value = '  +8.70%  '
$ws.Range("E51").Value = '  +1.27%  '
